# Updated cryptos list (price + volume refresh, with a few rows re-ranked)
# Note: Price cells in column D are stored as text in the source sheet
# (e.g. "232.80", "0.0240") to preserve exact formatting/trailing zeros,
# so their literals below are prefixed with a leading apostrophe
# (PowerShell '' -> literal ') which Excel interprets as "force text",
# matching the original text-typed cells instead of auto-converting to a
# number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''41.860.36'
$ws.Range("E2").Value = '  +1.73%  '

$ws.Range("D3").Value = '''2.229.26'
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").Value = '''232.80'
$ws.Range("E5").Value = '  +1.45%  '

$ws.Range("E6").Value = '  -1.68%  '

$ws.Range("D7").Value = '''60.42'
$ws.Range("E7").Value = '  -7.15%  '

$ws.Range("E8").Value = '  +0.20%  '

$ws.Range("D9").Value = '''0.404'
$ws.Range("E9").Value = '  -1.10%  '

$ws.Range("D10").Value = '''58.24'
$ws.Range("E10").Value = '  -3.33%  '

$ws.Range("E11").Value = '  +1.57%  '

$ws.Range("E12").Value = '  -0.36%  '

$ws.Range("D13").Value = '''2.559.17'
$ws.Range("E13").Value = '  +0.16%  '

$ws.Range("D14").Value = '''15.53'
$ws.Range("E14").Value = '  -3.77%  '

$ws.Range("D15").Value = '''22.79'
$ws.Range("E15").Value = '  +1.74%  '

$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '''5.63'
$ws.Range("E16").Value = '  -0.39%  '

$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").Value = '''0.802'
$ws.Range("E17").Value = '  -3.43%  '

$ws.Range("D18").Value = '''2.232.77'
$ws.Range("E18").Value = '  -0.06%  '

$ws.Range("D19").Value = '''41.678.20'
$ws.Range("E19").Value = '  +1.89%  '

$ws.Range("E20").Value = '  -0.19%  '

$ws.Range("D21").Value = '''72.51'
$ws.Range("E21").Value = '  -2.20%  '

$ws.Range("E22").Value = '  -0.29%  '

$ws.Range("D23").Value = '''247.92'
$ws.Range("E23").Value = '  -3.50%  '

$ws.Range("E24").Value = '  -0.20%  '

$ws.Range("E25").Value = '  +0.23%  '

$ws.Range("E26").Value = '  -3.28%  '

$ws.Range("D27").Value = '''9.65'
$ws.Range("E27").Value = '  -0.80%  '

$ws.Range("D28").Value = '''169.29'
$ws.Range("E28").Value = '  -2.30%  '

$ws.Range("E29").Value = '  -3.40%  '

$ws.Range("D30").Value = '''19.93'
$ws.Range("E30").Value = '  -2.55%  '

$ws.Range("E31").Value = '  -3.16%  '

$ws.Range("D32").Value = '''2.59'
$ws.Range("E32").Value = '  -8.64%  '

$ws.Range("E33").Value = '  -1.90%  '

$ws.Range("D34").Value = '''5.02'
$ws.Range("E34").Value = '  +3.13%  '

$ws.Range("E35").Value = '  +0.35%  '

$ws.Range("D36").Value = '''0.0657'
$ws.Range("E36").Value = '  +3.50%  '

$ws.Range("E37").Value = '  -9.17%  '

$ws.Range("D38").Value = '''2.39'
$ws.Range("E38").Value = '  -3.41%  '

$ws.Range("D39").Value = '''3.62'
$ws.Range("E39").Value = '  -5.46%  '

$ws.Range("E40").Value = '  +0.65%  '

$ws.Range("B41").Value = 'TerraClassic'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D41").Value = '''0.000235'
$ws.Range("E41").Value = '  -0.61%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.0240'
$ws.Range("E42").Value = '  +1.34%  '

$ws.Range("D43").Value = '''8.57'
$ws.Range("E43").Value = '  -2.62%  '

$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").Value = '''4.54'
$ws.Range("E44").Value = '  -6.56%  '

$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '''1.23'
$ws.Range("E45").Value = '  -1.63%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '''0.0968'
$ws.Range("E46").Value = '  +2.72%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '''98.89'
$ws.Range("E47").Value = '  -3.28%  '

$ws.Range("D48").Value = '''1.468.53'
$ws.Range("E48").Value = '  -2.90%  '

$ws.Range("D49").Value = '''16.63'
$ws.Range("E49").Value = '  -5.28%  '

$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").Value = '''2.80'
$ws.Range("E50").Value = '  -1.35%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '''2.29'
$ws.Range("E51").Value = '  +8.14%  '
